$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the test case name (B5) from "open URL to online store"
# to "open Home page to Online store"
$ws.Range("B5").Value = "open Home page to Online store"

# Move the active selection to B14, matching the saved selection state
$ws.Range("B14").Select() | Out-Null
